# Corrección a Diebold Mariano: actualizar valores de estadístico (col C) y p-value (col D)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.6712789683070062
$ws.Range("D2").Value = 0.5090304226709601

$ws.Range("C3").Value = 0.7408175863929127
$ws.Range("D3").Value = 0.4666381540943689

$ws.Range("C4").Value = 0.8101593474220016
$ws.Range("D4").Value = 0.4265256190582978

$ws.Range("C5").Value = -0.04449207906000651
$ws.Range("D5").Value = 0.9649135426589905

$ws.Range("C6").Value = 0.0870727984777056
$ws.Range("D6").Value = 0.9314013241124666

$ws.Range("C7").Value = 0.2545483706006483
$ws.Range("D7").Value = 0.8014361943774171

$ws.Range("C8").Value = -0.5317481000201038
$ws.Range("D8").Value = 0.6002294224763838

$ws.Range("C9").Value = 0.2152225572789606
$ws.Range("D9").Value = 0.8315765101758759

$ws.Range("C10").Value = -0.6597251512177
$ws.Range("D10").Value = 0.5162778690317205

$ws.Range("C11").Value = -0.7853407250842371
$ws.Range("D11").Value = 0.4406300627594613

$wb.Save()
